# Update worksheet values to reflect new TPM-derived statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending=FAPs, Target=ECs)
$ws.Cells.Item(2, 13).Value = 7.413580666666667   # M2
$ws.Cells.Item(2, 14).Value = 22.240742           # N2
$ws.Cells.Item(2, 15).Value = 0.05108888817597561 # O2
$ws.Cells.Item(2, 16).Value = 0.05108888817597561 # P2
$ws.Cells.Item(2, 17).Value = 31.31893347285022   # Q2
$ws.Cells.Item(2, 18).Value = 281.870401255652    # R2
$ws.Cells.Item(2, 19).Value = 0.03598305592682585 # S2
$ws.Cells.Item(2, 20).Value = 0.03598305592682585 # T2

# Row 3 (Sending=FAPs, Target=FAPs)
$ws.Cells.Item(3, 15).Value = 0.0112127179963522   # O3
$ws.Cells.Item(3, 16).Value = 0.0112127179963522   # P3
$ws.Cells.Item(3, 19).Value = 0.007897370116271141 # S3
$ws.Cells.Item(3, 20).Value = 0.007897370116271143 # T3

# Row 4 (Sending=FAPs, Target=MuSCs)
$ws.Cells.Item(4, 15).Value = 0.9376983938276722  # O4
$ws.Cells.Item(4, 16).Value = 0.9376983938276722  # P4
$ws.Cells.Item(4, 19).Value = 0.6604421225878745  # S4
$ws.Cells.Item(4, 20).Value = 0.6604421225878745  # T4

# Row 5 (Sending=MuSCs, Target=ECs)
$ws.Cells.Item(5, 13).Value = 7.413580666666667   # M5
$ws.Cells.Item(5, 14).Value = 22.240742           # N5
$ws.Cells.Item(5, 15).Value = 0.05108888817597561 # O5
$ws.Cells.Item(5, 16).Value = 0.05108888817597561 # P5
$ws.Cells.Item(5, 17).Value = 13.147814799978     # Q5
$ws.Cells.Item(5, 18).Value = 118.330333199802    # R5
$ws.Cells.Item(5, 19).Value = 0.01510583224914977 # S5
$ws.Cells.Item(5, 20).Value = 0.01510583224914976 # T5

# Row 6 (Sending=MuSCs, Target=FAPs)
$ws.Cells.Item(6, 15).Value = 0.0112127179963522   # O6
$ws.Cells.Item(6, 16).Value = 0.0112127179963522   # P6
$ws.Cells.Item(6, 19).Value = 0.003315347880081059 # S6
$ws.Cells.Item(6, 20).Value = 0.003315347880081059 # T6

# Row 7 (Sending=MuSCs, Target=MuSCs)
$ws.Cells.Item(7, 15).Value = 0.9376983938276722  # O7
$ws.Cells.Item(7, 16).Value = 0.9376983938276722  # P7
$ws.Cells.Item(7, 19).Value = 0.2772562712397977  # S7
$ws.Cells.Item(7, 20).Value = 0.2772562712397977  # T7
